$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "all": insert new row 45 (copy of row 44, next day, no new cases)
# ----------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all")
$null = $wsAll.Rows.Item(45).Insert()

$wsAll.Range("A45").Value = 43973
$wsAll.Range("B45").Value = 285
$wsAll.Range("C45").Value = 282
$wsAll.Range("D45").Value = 37
$wsAll.Range("E45").Value = 32
$wsAll.Range("F45").Value = 5
$wsAll.Range("G45").Value = 11
$wsAll.Range("H45").Value = 234

# ----------------------------------------------------------------------
# Sheet "kobe": insert new row 100 (next day; only new-case count changes)
# ----------------------------------------------------------------------
$wsKobe = $wb.Worksheets.Item("kobe")
$null = $wsKobe.Rows.Item(100).Insert()

$wsKobe.Range("A100").Value = 43973
$wsKobe.Range("B100").Value = 0
$wsKobe.Range("C100").Value = 2982
$wsKobe.Range("D100").Value = 0
$wsKobe.Range("E100").Value = 285
$wsKobe.Range("F100").Value = 32
$wsKobe.Range("G100").Value = 28
$wsKobe.Range("H100").Value = 4
$wsKobe.Range("I100").Value = 11
$wsKobe.Range("J100").Value = 225

# ----------------------------------------------------------------------
# Sheet "other": insert new row 75 (identical totals to row 74, next day)
# ----------------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("other")
$null = $wsOther.Rows.Item(75).Insert()

$wsOther.Range("A75").Value = 43973
$wsOther.Range("B75").Value = 0
$wsOther.Range("C75").Value = 14
$wsOther.Range("D75").Value = 5
$wsOther.Range("E75").Value = 4
$wsOther.Range("F75").Value = 1
$wsOther.Range("G75").Value = 0
$wsOther.Range("H75").Value = 9

# ----------------------------------------------------------------------
# View state: update each sheet's selection to rest on the new last
# data row, and make "all" the active sheet/tab (as in the target file).
# ----------------------------------------------------------------------
$null = $wsKobe.Range("B101").Select()
$null = $wsOther.Range("B76").Select()
$null = $wsAll.Activate()
$null = $wsAll.Range("J45").Select()
